# Apply the "fixed utilisations and costs" edit:
#  - rename the laser_1..4 / surgery_1..4 / surgery_5 item codes in column A
#    (rows 2-10) to laser_treatment_1..4 / surgery_treatment_1..5
#  - normalise the banding/border format of column A for rows 2-10 so that
#    all "laser_*" rows (2-5) share one look and all "surgery_*" rows (6-10)
#    share another
#  - grow row 10 to fit the longer "surgery_treatment_5" label
#  - move the active selection to A11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A item-code renames (rows 2-10) ---------------------------------
$ws.Range("A2").Value  = "laser_treatment_1"
$ws.Range("A3").Value  = "laser_treatment_2"
$ws.Range("A4").Value  = "laser_treatment_3"
$ws.Range("A5").Value  = "laser_treatment_4"
$ws.Range("A6").Value  = "surgery_treatment_1"
$ws.Range("A7").Value  = "surgery_treatment_2"
$ws.Range("A8").Value  = "surgery_treatment_3"
$ws.Range("A9").Value  = "surgery_treatment_4"
$ws.Range("A10").Value = "surgery_treatment_5"

# --- Re-align the direct formatting of column A ------------------------------
# Rows 2-5 ("laser_*") all take on the banding used by A2.
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats

# Rows 6-10 ("surgery_*") all take on the banding used by A6.
$ws.Range("A6").Copy()
$ws.Range("A7:A10").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Row height fix for the now-longer row 10 label --------------------------
$ws.Rows.Item(10).RowHeight = 36

# --- Move the selection, matching the saved cursor position ------------------
$ws.Range("A11").Select()
